$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.963.52'
$ws.Range("E2").Value = '  -1.43%  '

# Row 3
$ws.Range("D3").Value = '1.818.56'
$ws.Range("E3").Value = '  -0.78%  '

# Row 4
$ws.Range("E4").Value = '  -0.23%  '

# Row 5
$ws.Range("E5").Value = '  -0.12%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '308.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.73%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4658'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.73%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3652'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.13%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07222'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.17%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8591'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.05%  '

# Row 11
$ws.Range("E11").Value = '  -3.46%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07562'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.18%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.846.07'
$ws.Range("E13").Value = '  -3.43%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.321'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.42%  '

# Row 15
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.44%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.461'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.91%  '

# Row 17
$ws.Range("E17").Value = '  +0.09%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008615'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.24%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.007'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.27%  '

# Row 21
$ws.Range("D21").Value = '26.566.15'
$ws.Range("E21").Value = '  -3.61%  '

# Row 22
$ws.Range("E22").Value = '  -3.29%  '

# Row 23
$ws.Range("E23").Value = '  -1.47%  '

# Row 24
$ws.Range("D24").Value = '1.998.10'
$ws.Range("E24").Value = '  -5.60%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.848'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.84%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.99%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.065'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.64%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.096'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.82%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.95%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08868'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.46%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.961'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.64%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.416'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.96%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.129'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.31%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7156'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.31%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.076'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.58%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05249'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.96%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01921'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.71%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.919'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.08%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.373'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.03%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.136'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.51%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5144'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.34%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1622'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.36%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.148'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.07%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4808'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.14%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.009'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.13%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.78%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.23%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06249'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.75%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.615'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.63%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.50%  '
